$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.546.94'
Set-TextValue $ws.Range('E2') '  -0.04%  '
Set-TextValue $ws.Range('D3') '1.753.62'
Set-TextValue $ws.Range('E3') '  +0.09%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '324.27'
Set-TextValue $ws.Range('E5') '  +0.12%  '
Set-TextValue $ws.Range('E6') '  -0.05%  '
Set-TextValue $ws.Range('D7') '0.4576'
Set-TextValue $ws.Range('E7') '  +2.93%  '
Set-TextValue $ws.Range('D8') '0.3560'
Set-TextValue $ws.Range('E8') '  -1.36%  '
Set-TextValue $ws.Range('D9') '0.07470'
Set-TextValue $ws.Range('E9') '  -0.19%  '
Set-TextValue $ws.Range('D10') '41.47'
Set-TextValue $ws.Range('E10') '  -1.78%  '
Set-TextValue $ws.Range('D11') '1.085'
Set-TextValue $ws.Range('E11') '  -1.82%  '
Set-TextValue $ws.Range('E12') '  -0.07%  '
Set-TextValue $ws.Range('D13') '20.75'
Set-TextValue $ws.Range('E13') '  +0.41%  '
Set-TextValue $ws.Range('D14') '6.007'
Set-TextValue $ws.Range('E14') '  -0.39%  '
Set-TextValue $ws.Range('D15') '7.172'
Set-TextValue $ws.Range('E15') '  -0.10%  '
Set-TextValue $ws.Range('D16') '1.756.96'
Set-TextValue $ws.Range('E16') '  -0.01%  '
Set-TextValue $ws.Range('D17') '94.23'
Set-TextValue $ws.Range('E17') '  +1.38%  '
Set-TextValue $ws.Range('D18') '0.00001055'
Set-TextValue $ws.Range('E18') '  -0.73%  '
Set-TextValue $ws.Range('D19') '0.06414'
Set-TextValue $ws.Range('E19') '  -0.07%  '
Set-TextValue $ws.Range('D20') '1.000'
Set-TextValue $ws.Range('E20') '  -0.10%  '
Set-TextValue $ws.Range('D21') '17.09'
Set-TextValue $ws.Range('E21') '  +0.36%  '
Set-TextValue $ws.Range('D22') '5.740'
Set-TextValue $ws.Range('E22') '  -1.68%  '
Set-TextValue $ws.Range('D23') '27.601.12'
Set-TextValue $ws.Range('D24') '11.19'
Set-TextValue $ws.Range('E24') '  -0.51%  '
Set-TextValue $ws.Range('D25') '2.084'
Set-TextValue $ws.Range('E25') '  -0.86%  '
Set-TextValue $ws.Range('D26') '165.50'
Set-TextValue $ws.Range('E26') '  +1.57%  '
Set-TextValue $ws.Range('D27') '20.15'
Set-TextValue $ws.Range('E27') '  -1.26%  '
Set-TextValue $ws.Range('D28') '1.953.88'
Set-TextValue $ws.Range('E28') '  -0.10%  '
Set-TextValue $ws.Range('D29') '2.118'
Set-TextValue $ws.Range('E29') '  -0.33%  '
Set-TextValue $ws.Range('D30') '125.66'
Set-TextValue $ws.Range('E30') '  +0.07%  '
Set-TextValue $ws.Range('D31') '1.080'
Set-TextValue $ws.Range('E31') '  -0.34%  '
Set-TextValue $ws.Range('D32') '0.09227'
Set-TextValue $ws.Range('E32') '  +2.30%  '
Set-TextValue $ws.Range('D33') '3.656'
Set-TextValue $ws.Range('E33') '  +0.51%  '
Set-TextValue $ws.Range('D34') '5.516'
Set-TextValue $ws.Range('E34') '  -0.45%  '
Set-TextValue $ws.Range('D35') '0.02284'
Set-TextValue $ws.Range('E35') '  -0.75%  '
Set-TextValue $ws.Range('D36') '11.71'
Set-TextValue $ws.Range('E36') '  -3.10%  '
Set-TextValue $ws.Range('D37') '0.2089'
Set-TextValue $ws.Range('E37') '  -0.30%  '
Set-TextValue $ws.Range('D38') '0.06014'
Set-TextValue $ws.Range('E38') '  +1.01%  '
Set-TextValue $ws.Range('D39') '0.6287'
Set-TextValue $ws.Range('E39') '  -1.02%  '
Set-TextValue $ws.Range('D40') '4.921'
Set-TextValue $ws.Range('E40') '  -0.27%  '
Set-TextValue $ws.Range('D41') '1.180'
Set-TextValue $ws.Range('E41') '  -1.07%  '
Set-TextValue $ws.Range('E42') '  -0.01%  '
Set-TextValue $ws.Range('D43') '7.792'
Set-TextValue $ws.Range('E43') '  +0.00%  '
Set-TextValue $ws.Range('D44') '13.13'
Set-TextValue $ws.Range('E44') '  -0.54%  '
Set-TextValue $ws.Range('D45') '3.717'
Set-TextValue $ws.Range('E45') '  +0.08%  '
Set-TextValue $ws.Range('D46') '0.5852'
Set-TextValue $ws.Range('E46') '  -0.26%  '
Set-TextValue $ws.Range('D47') '121.94'
Set-TextValue $ws.Range('E47') '  +0.33%  '
Set-TextValue $ws.Range('D48') '1.935'
Set-TextValue $ws.Range('E48') '  -1.31%  '
Set-TextValue $ws.Range('D49') '0.06890'
Set-TextValue $ws.Range('E49') '  +0.59%  '
Set-TextValue $ws.Range('E50') '  -2.46%  '
Set-TextValue $ws.Range('D51') '71.92'
Set-TextValue $ws.Range('E51') '  -0.64%  '
